# Auto-generated edit script: update Leve market-price snapshot columns (H:N)
# across ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled-runner commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 276.8
$ws.Range("I2").Value = 84
$ws.Range("J2").Value = 662.4
$ws.Range("K2").Value = 84
$ws.Range("L2").Value = 662.4
$ws.Range("M2").Value = 29
$ws.Range("N2").Value = -888.4

$ws.Range("H4").Value = 398
$ws.Range("I4").Value = 398
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 398
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -284
$ws.Range("N4").Value = ""

$ws.Range("H7").Value = 13529.417
$ws.Range("I7").Value = 1152.5
$ws.Range("J7").Value = 16004.8
$ws.Range("K7").Value = 1152.5
$ws.Range("L7").Value = 16004.8
$ws.Range("M7").Value = -1040.5
$ws.Range("N7").Value = -16228.8

$ws.Range("H8").Value = 20
$ws.Range("I8").Value = 20
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 60
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 79
$ws.Range("N8").Value = ""

$ws.Range("H9").Value = 413
$ws.Range("I9").Value = 800
$ws.Range("J9").Value = 284
$ws.Range("K9").Value = 800
$ws.Range("L9").Value = 284
$ws.Range("M9").Value = -631
$ws.Range("N9").Value = -622

$ws.Range("H10").Value = 10000
$ws.Range("J10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("N10").Value = -10586

$ws.Range("H14").Value = 13529.417
$ws.Range("I14").Value = 1152.5
$ws.Range("J14").Value = 16004.8
$ws.Range("K14").Value = 1152.5
$ws.Range("L14").Value = 16004.8
$ws.Range("M14").Value = -961.5
$ws.Range("N14").Value = -16386.8

$ws.Range("H16").Value = 16596
$ws.Range("I16").Value = 4980
$ws.Range("J16").Value = 19500
$ws.Range("K16").Value = 4980
$ws.Range("L16").Value = 19500
$ws.Range("M16").Value = -4750
$ws.Range("N16").Value = -19960

$ws.Range("H18").Value = 954.7273
$ws.Range("I18").Value = 900
$ws.Range("J18").Value = 1201
$ws.Range("K18").Value = 900
$ws.Range("L18").Value = 1201
$ws.Range("M18").Value = -616
$ws.Range("N18").Value = -1769

$ws.Range("H21").Value = 34253.168
$ws.Range("I21").Value = 44759.5
$ws.Range("J21").Value = 29000
$ws.Range("K21").Value = 44759.5
$ws.Range("L21").Value = 29000
$ws.Range("M21").Value = -44291.5
$ws.Range("N21").Value = -29936

$ws.Range("H23").Value = 34253.168
$ws.Range("I23").Value = 44759.5
$ws.Range("J23").Value = 29000
$ws.Range("K23").Value = 44759.5
$ws.Range("L23").Value = 29000
$ws.Range("M23").Value = -44525.5
$ws.Range("N23").Value = -29468

$ws.Range("H32").Value = 1012.3077
$ws.Range("I32").Value = 935
$ws.Range("J32").Value = 1136
$ws.Range("K32").Value = 935
$ws.Range("L32").Value = 1136
$ws.Range("M32").Value = -609
$ws.Range("N32").Value = -1788

$ws.Range("I33").Value = 228.64706
$ws.Range("K33").Value = 228.64706
$ws.Range("M33").Value = 0.3529399999999896

$ws.Range("H38").Value = 49.4
$ws.Range("I38").Value = 49.4
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 148.2
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 223.8
$ws.Range("N38").Value = ""

$ws.Range("H39").Value = 832.2222
$ws.Range("I39").Value = 126.888885
$ws.Range("J39").Value = 1537.5555
$ws.Range("K39").Value = 380.666655
$ws.Range("L39").Value = 4612.666499999999
$ws.Range("M39").Value = -84.66665499999999
$ws.Range("N39").Value = -5204.666499999999

$ws.Range("H40").Value = 1649.1111
$ws.Range("I40").Value = 1380.64
$ws.Range("J40").Value = 2259.2727
$ws.Range("K40").Value = 1380.64
$ws.Range("L40").Value = 2259.2727
$ws.Range("M40").Value = -1205.64
$ws.Range("N40").Value = -2609.2727

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = ""

$ws.Range("H48").Value = 2519
$ws.Range("J48").Value = 2519
$ws.Range("L48").Value = 7557
$ws.Range("N48").Value = -8141

$ws.Range("H51").Value = 5200
$ws.Range("I51").Value = 2425
$ws.Range("J51").Value = 5992.857
$ws.Range("K51").Value = 2425
$ws.Range("L51").Value = 5992.857
$ws.Range("M51").Value = -1941
$ws.Range("N51").Value = -6960.857

$ws.Range("H56").Value = 2519
$ws.Range("J56").Value = 2519
$ws.Range("L56").Value = 7557
$ws.Range("N56").Value = -8625

$ws.Range("H58").Value = 2005.3636
$ws.Range("I58").Value = 2165.5715
$ws.Range("J58").Value = 1725
$ws.Range("K58").Value = 6496.7145
$ws.Range("L58").Value = 5175
$ws.Range("M58").Value = -6346.7145
$ws.Range("N58").Value = -5475

$ws.Range("H100").Value = 3182.5
$ws.Range("I100").Value = 3091.6667
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 3091.6667
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -2550.6667
$ws.Range("N100").Value = -5082

$ws.Range("H116").Value = 2719.4443
$ws.Range("I116").Value = 2470
$ws.Range("J116").Value = 2878.182
$ws.Range("K116").Value = 2470
$ws.Range("L116").Value = 2878.182
$ws.Range("M116").Value = 972
$ws.Range("N116").Value = -9762.182000000001

$ws.Range("H129").Value = 634.3125
$ws.Range("I129").Value = 407.5
$ws.Range("K129").Value = 1222.5
$ws.Range("M129").Value = 3777.5

$ws.Range("H137").Value = 2018.4
$ws.Range("I137").Value = 1499.7142
$ws.Range("J137").Value = 2297.6924
$ws.Range("K137").Value = 4499.142599999999
$ws.Range("L137").Value = 6893.0772
$ws.Range("M137").Value = -1949.142599999999
$ws.Range("N137").Value = -11993.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2004
$ws.Range("I61").Value = 1931.5555
$ws.Range("K61").Value = 1931.5555
$ws.Range("M61").Value = -1719.5555

$ws.Range("H133").Value = 33799.8
$ws.Range("J133").Value = 33799.8
$ws.Range("L133").Value = 33799.8
$ws.Range("N133").Value = -38859.8

$ws.Range("H136").Value = 2004
$ws.Range("I136").Value = 1931.5555
$ws.Range("K136").Value = 5794.666499999999
$ws.Range("M136").Value = -3244.666499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 2000
$ws.Range("J29").Value = 2000
$ws.Range("L29").Value = 2000
$ws.Range("N29").Value = -2586

$ws.Range("H86").Value = 2475956
$ws.Range("I86").Value = 4836748.5
$ws.Range("J86").Value = 7854.8184
$ws.Range("K86").Value = 4836748.5
$ws.Range("L86").Value = 7854.8184
$ws.Range("M86").Value = -4835625.5
$ws.Range("N86").Value = -10100.8184

$ws.Range("H89").Value = 2475956
$ws.Range("I89").Value = 4836748.5
$ws.Range("J89").Value = 7854.8184
$ws.Range("K89").Value = 24183742.5
$ws.Range("L89").Value = 39274.092
$ws.Range("M89").Value = -24178126.5
$ws.Range("N89").Value = -50506.092

$ws.Range("H140").Value = 42352.223
$ws.Range("J140").Value = 42352.223
$ws.Range("L140").Value = 42352.223
$ws.Range("N140").Value = -52712.223

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 865.47
$ws.Range("J131").Value = 891.663
$ws.Range("L131").Value = 2674.989
$ws.Range("N131").Value = -12754.989

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1184.8235
$ws.Range("I107").Value = 826.8889
$ws.Range("K107").Value = 826.8889
$ws.Range("M107").Value = 1093.1111

$ws.Range("H126").Value = 45459810
$ws.Range("I126").Value = 66672400
$ws.Range("J126").Value = 4244.857
$ws.Range("K126").Value = 200017200
$ws.Range("L126").Value = 12734.571
$ws.Range("M126").Value = -200014730
$ws.Range("N126").Value = -17674.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 738.52
$ws.Range("I22").Value = 588.7
$ws.Range("J22").Value = 838.4
$ws.Range("K22").Value = 588.7
$ws.Range("L22").Value = 838.4
$ws.Range("M22").Value = -293.7
$ws.Range("N22").Value = -1428.4

$ws.Range("H27").Value = 738.52
$ws.Range("I27").Value = 588.7
$ws.Range("J27").Value = 838.4
$ws.Range("K27").Value = 588.7
$ws.Range("L27").Value = 838.4
$ws.Range("M27").Value = -481.7
$ws.Range("N27").Value = -1052.4

$ws.Range("H132").Value = 3020.4075
$ws.Range("I132").Value = 2550
$ws.Range("J132").Value = 4666.8335
$ws.Range("K132").Value = 7650
$ws.Range("L132").Value = 14000.5005
$ws.Range("M132").Value = -5120
$ws.Range("N132").Value = -19060.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 12008.333
$ws.Range("I40").Value = 6025
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 6025
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = -5876
$ws.Range("N40").Value = -15298

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = ""

$ws.Range("H50").Value = 12000
$ws.Range("J50").Value = 12000
$ws.Range("L50").Value = 12000
$ws.Range("N50").Value = -13262

$ws.Range("H51").Value = 10000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 10000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 10000
$ws.Range("M51").Value = ""
$ws.Range("N51").Value = -11020

$ws.Range("H52").Value = 8542.143
$ws.Range("I52").Value = 4500
$ws.Range("J52").Value = 9215.833000000001
$ws.Range("K52").Value = 4500
$ws.Range("L52").Value = 9215.833000000001
$ws.Range("M52").Value = -4274
$ws.Range("N52").Value = -9667.833000000001

$ws.Range("H58").Value = 7932.6665
$ws.Range("I58").Value = 4000
$ws.Range("J58").Value = 9899
$ws.Range("K58").Value = 4000
$ws.Range("L58").Value = 9899
$ws.Range("M58").Value = -3692
$ws.Range("N58").Value = -10515

$ws.Range("H132").Value = 3203.543
$ws.Range("I132").Value = 3192.4348
$ws.Range("J132").Value = 3224.8333
$ws.Range("K132").Value = 9577.304400000001
$ws.Range("L132").Value = 9674.499899999999
$ws.Range("M132").Value = -7047.304400000001
$ws.Range("N132").Value = -14734.4999
